$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1976047904191617
$ws.Range("C2").Value = 0.5209580838323353
$ws.Range("J2").Value = 0.01796407185628742
$ws.Range("P2").Value = 0.1287425149700599
$ws.Range("S2").Value = 0.1347305389221557
$ws.Range("B3").Value = 0.01098901098901099
$ws.Range("C3").Value = 0.03846153846153846
$ws.Range("J3").Value = 0.01648351648351648
$ws.Range("P3").Value = 0.7252747252747253
$ws.Range("S3").Value = 0.2087912087912088
$ws.Range("B6").Value = 0.04942965779467681
$ws.Range("D6").Value = 0.007604562737642586
$ws.Range("F6").Value = 0.05703422053231939
$ws.Range("J6").Value = 0.2167300380228137
$ws.Range("O6").Value = 0.03422053231939164
$ws.Range("Q6").Value = 0.1673003802281369
$ws.Range("R6").Value = 0.07224334600760456
$ws.Range("S6").Value = 0.3954372623574144
$ws.Range("B7").Value = 0.1222707423580786
$ws.Range("D7").Value = 0.01746724890829694
$ws.Range("F7").Value = 0.07423580786026202
$ws.Range("J7").Value = 0.1703056768558952
$ws.Range("O7").Value = 0.02620087336244541
$ws.Range("Q7").Value = 0.148471615720524
$ws.Range("R7").Value = 0.09606986899563319
$ws.Range("S7").Value = 0.3449781659388647
$ws.Range("B8").Value = 0.1054852320675106
$ws.Range("D8").Value = 0.02320675105485232
$ws.Range("E8").Value = 0.002109704641350211
$ws.Range("F8").Value = 0.08016877637130802
$ws.Range("J8").Value = 0.120253164556962
$ws.Range("O8").Value = 0.01054852320675105
$ws.Range("Q8").Value = 0.1962025316455696
$ws.Range("R8").Value = 0.0759493670886076
$ws.Range("S8").Value = 0.3860759493670886
$ws.Range("B9").Value = 0.08035714285714286
$ws.Range("D9").Value = 0.01339285714285714
$ws.Range("E9").Value = 0.004464285714285714
$ws.Range("F9").Value = 0.07589285714285714
$ws.Range("J9").Value = 0.1294642857142857
$ws.Range("O9").Value = 0.04017857142857143
$ws.Range("Q9").Value = 0.1651785714285714
$ws.Range("R9").Value = 0.1071428571428571
$ws.Range("S9").Value = 0.3839285714285715
$ws.Range("B10").Value = 0.1034246575342466
$ws.Range("D10").Value = 0.02191780821917808
$ws.Range("F10").Value = 0.07671232876712329
$ws.Range("J10").Value = 0.1047945205479452
$ws.Range("O10").Value = 0.02534246575342466
$ws.Range("Q10").Value = 0.1958904109589041
$ws.Range("R10").Value = 0.08904109589041095
$ws.Range("S10").Value = 0.3828767123287671
$ws.Range("G11").Value = 0.1613756613756614
$ws.Range("J11").Value = 0.1084656084656085
$ws.Range("K11").Value = 0.2169312169312169
$ws.Range("L11").Value = 0.4947089947089947
$ws.Range("S11").Value = 0.01851851851851852
$ws.Range("G12").Value = 0.7433155080213903
$ws.Range("J12").Value = 0.2299465240641711
$ws.Range("K12").Value = 0.0053475935828877
$ws.Range("S12").Value = 0.0213903743315508
$ws.Range("G13").Value = 0.6428571428571429
$ws.Range("J13").Value = 0.2678571428571428
$ws.Range("S13").Value = 0.08928571428571429
$ws.Range("F15").Value = 0.01310043668122271
$ws.Range("H15").Value = 0.1528384279475982
$ws.Range("I15").Value = 0.03493449781659388
$ws.Range("J15").Value = 0.3275109170305677
$ws.Range("K15").Value = 0.07860262008733625
$ws.Range("M15").Value = 0.02620087336244541
$ws.Range("O15").Value = 0.03930131004366812
$ws.Range("S15").Value = 0.3275109170305677
$ws.Range("F16").Value = 0.02450980392156863
$ws.Range("H16").Value = 0.1470588235294118
$ws.Range("I16").Value = 0.06862745098039216
$ws.Range("J16").Value = 0.4509803921568628
$ws.Range("K16").Value = 0.1225490196078431
$ws.Range("M16").Value = 0.02941176470588235
$ws.Range("O16").Value = 0.01470588235294118
$ws.Range("S16").Value = 0.142156862745098
$ws.Range("F17").Value = 0.01431492842535787
$ws.Range("H17").Value = 0.1574642126789366
$ws.Range("I17").Value = 0.1104294478527607
$ws.Range("J17").Value = 0.4335378323108384
$ws.Range("K17").Value = 0.08793456032719836
$ws.Range("M17").Value = 0.02862985685071575
$ws.Range("N17").Value = 0.002044989775051125
$ws.Range("O17").Value = 0.049079754601227
$ws.Range("S17").Value = 0.1165644171779141
$ws.Range("F18").Value = 0.01298701298701299
$ws.Range("H18").Value = 0.2121212121212121
$ws.Range("I18").Value = 0.1212121212121212
$ws.Range("J18").Value = 0.4025974025974026
$ws.Range("K18").Value = 0.09090909090909091
$ws.Range("M18").Value = 0.02164502164502164
$ws.Range("O18").Value = 0.0303030303030303
$ws.Range("S18").Value = 0.1082251082251082
$ws.Range("F19").Value = 0.01131824234354194
$ws.Range("H19").Value = 0.1917443408788282
$ws.Range("I19").Value = 0.08189081225033289
$ws.Range("J19").Value = 0.374167776298269
$ws.Range("K19").Value = 0.1191744340878828
$ws.Range("M19").Value = 0.01930758988015979
$ws.Range("O19").Value = 0.05858854860186418
$ws.Range("S19").Value = 0.1438082556591212
